$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Merchant (商店) - add repeat limit to the "discard item for gold" effect
$ws.Range("C8").Value = "多选：①可重复至多3次，弃置1张物品牌，然后获得1金币。②取遭遇牌堆前3张物品牌，然后玩家每支付2金币，可以获得其中1张牌。"

# Row 9: Trainer (训练场) - turn into a multi-choice with a new reroll effect
$ws.Range("C9").Value = "多选：①可重复，支付3金币，从购买能力区选1张牌获得。②可重复至多3次，支付1金币，将购买能力区补满，然后重抽其中任意张牌。"

# Row 9 now wraps to two lines, so double its height to match
$ws.Range("A9:F9").RowHeight = 57

# Move the active selection to E9 (matches the saved cursor position in the edit)
$ws.Range("E9").Select()
